$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 310 (E310, F310 changed)
$ws.Range("E310").Value = 3.7041
$ws.Range("F310").Value = 3.7051

# Add new row 311
$ws.Range("A311").Value = 45047.33333333334
$ws.Range("B311").Value = "FX_IDC:USDPEN"
$ws.Range("C311").Value = 3.7026
$ws.Range("D311").Value = 3.7503
$ws.Range("E311").Value = 3.6386
$ws.Range("F311").Value = 3.6559
$ws.Range("G311").Value = 0

# Add new row 312
$ws.Range("A312").Value = 45078.33333333334
$ws.Range("B312").Value = "FX_IDC:USDPEN"
$ws.Range("C312").Value = 3.6559
$ws.Range("D312").Value = 3.6947
$ws.Range("E312").Value = 3.6008
$ws.Range("F312").Value = 3.61
$ws.Range("G312").Value = 0

# Add new row 313
$ws.Range("A313").Value = 45110.33333333334
$ws.Range("B313").Value = "FX_IDC:USDPEN"
$ws.Range("C313").Value = 3.6089
$ws.Range("D313").Value = 3.6488
$ws.Range("E313").Value = 3.6069
$ws.Range("F313").Value = 3.6264
$ws.Range("G313").Value = 0

# Apply same style as column A date cells (row 310) to the new date cells
$ws.Range("A310").Copy() | Out-Null
$ws.Range("A311:A313").PasteSpecial(-4122) | Out-Null
